$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# A new item (BABETONE SYRUP SUGAR FREE 120 ML) was sold, so it needs to be
# inserted as the first row of the items list (row 7), pushing the existing
# two items (DEPO-PEN, TORSERETIC) down by one row, and the totals / footer
# rows down as well. The grand total and the footer timestamp are updated
# accordingly.
# ---------------------------------------------------------------------------

# Step 1: Insert a new blank row at row 7 (the top of the items list).
# This shifts DEPO-PEN (was row 7) -> row 8, TORSERETIC (was row 8) -> row 9,
# the totals row (was row 9) -> row 10, and the footer row (was row 10) -> row 11.
$ws.Rows("7:7").Insert()

# Step 2: Copy the formatting of an existing item row (row 8, DEPO-PEN) into
# the new row 7, cell by cell, so every column keeps the correct per-cell style.
$itemCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")
foreach ($col in $itemCols) {
    $ws.Range($col + "8").Copy()
    $ws.Range($col + "7").PasteSpecial(-4122)   # xlPasteFormats
}
$excel.CutCopyMode = $false

# Step 3: Re-create the merged cells for the new row 7 (matching the pattern
# used by the other item rows).
$ws.Range("A7:B7").Merge()
$ws.Range("C7:G7").Merge()
$ws.Range("H7:K7").Merge()
$ws.Range("L7:M7").Merge()
$ws.Range("N7:O7").Merge()

# Step 4: Fix up the row heights for the whole block.
$ws.Rows("7:7").RowHeight = 25.5
$ws.Rows("8:8").RowHeight = 24.75
$ws.Rows("9:9").RowHeight = 25.5
$ws.Rows("10:10").RowHeight = 24.75
$ws.Rows("11:11").RowHeight = 16.5

# Step 5: Fill in the data for the new first item row (BABETONE).
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "BABETONE SYRUP SUGAR FREE 120 ML"
$ws.Range("H7").Value = "0:0"
$ws.Range("L7").Value = "1"
$ws.Range("N7").Value = "35.00"
$ws.Range("P7").Value = "35.0000"
$ws.Range("Q7").Value = "1:0"

# Step 6: Renumber the items that shifted down (DEPO-PEN becomes item 2,
# TORSERETIC becomes item 3). Their other cell values are unchanged.
$ws.Range("A8").Value = 2
$ws.Range("A9").Value = 3

# Step 7: Update the grand total (was 111.13, now + 35.00 = 146.13).
$ws.Range("P10").Value = 146.13

# Step 8: Update the printed timestamp in the footer row.
$ws.Range("A11").Value = "Wednesday, 20 August, 2025 10:16 AM"
